$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Summary table, 10,000 Elements column (D)
$ws.Range("D4").Value = 5.8983800000000004
$ws.Range("D5").Value = 152.03659999999999
$ws.Range("D6").Value = 24.58466

# Hybrid (10,000 Elements) trial data, column H rows 10-14
$ws.Range("H10").Value = 7.5720999999999998
$ws.Range("H11").Value = 5.3691000000000004
$ws.Range("H12").Value = 5.5719000000000003
$ws.Range("H13").Value = 5.3932000000000002
$ws.Range("H14").Value = 5.5856000000000003

# Bubble (10,000 Elements) trial data, column H rows 17-21
$ws.Range("H17").Value = 150.71100000000001
$ws.Range("H18").Value = 151.33600000000001
$ws.Range("H19").Value = 150.72900000000001
$ws.Range("H20").Value = 158.339
$ws.Range("H21").Value = 149.06800000000001

# Merge (10,000 Elements) trial data, column H rows 24-28
$ws.Range("H24").Value = 25.076599999999999
$ws.Range("H25").Value = 24.748200000000001
$ws.Range("H26").Value = 24.3325
$ws.Range("H27").Value = 24.497900000000001
$ws.Range("H28").Value = 24.2681

# Update selection to match the post-edit state
[void]$ws.Range("I15").Select()
